# Commit: "added filter for the room type for get rooms"
#
# - Rename Sheet3 -> Enhancements, make it the active sheet.
# - Clear the old tabSelected/selection from Issues, select whole sheet there.
# - Populate the Enhancements sheet with the new "Enhancements" tracking table
#   (dates, page/url, description, assignee, found-on / start / end dates,
#   comments) describing the roomType filter work.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data fields")
$ws2 = $wb.Worksheets.Item("Issues")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------------
# 1. Rename the third sheet and rebuild its layout / data.
# ---------------------------------------------------------------------------
$ws3.Name = "Enhancements"

# Column widths matching the "Issues" sheet layout.
$ws2.Range("A1:H1").Copy()
$ws3.Range("A1:H1").PasteSpecial(-4104)   # xlPasteColumnWidths

# Header row (reuse existing shared strings from Issues' header row).
$ws2.Range("A1:H1").Copy()
$ws3.Range("A1:H1").PasteSpecial(-4122)   # xlPasteFormats (bring over the bold/shaded style)
$ws3.Range("A1").Value = "Date"
$ws3.Range("B1").Value = "Page / application Url"
$ws3.Range("C1").Value = "Description"
$ws3.Range("D1").Value = "Assignee"
$ws3.Range("E1").Value = "Found on"
$ws3.Range("F1").Value = "start date"
$ws3.Range("G1").Value = "endDate"
$ws3.Range("H1").Value = "Comments"

# Helper: bring over the date-cell number format (style 8 on Issues!A2) before
# typing, so the shared-string cells below keep a date style even when they
# hold plain text like "19-9-16" (mirrors how the source sheet was built).
function Set-DateStyle($range) {
    $ws2.Range("A2").Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
}

# Helper: bring over the wrapped-text style (style 1 on Issues!B2).
function Set-TextStyle($range) {
    $ws2.Range("B2").Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
}

# Row 2
Set-DateStyle($ws3.Range("A2"))
$ws3.Range("A2").Value = 42632
Set-TextStyle($ws3.Range("B2:D2"))
$ws3.Range("B2").Value = "/room(GET)"
$ws3.Range("C2").Value = "add extra parameter for the filter as roomType in the client side"
$ws3.Range("D2").Value = "rajashree"
Set-DateStyle($ws3.Range("E2:F2"))
$ws3.Range("E2").Value = "19-9-16"
$ws3.Range("F2").Value = "19-9-16"
$ws3.Rows.Item(2).RowHeight = 30

# Row 3
Set-DateStyle($ws3.Range("A3"))
$ws3.Range("A3").Value = 42632
Set-TextStyle($ws3.Range("B3:D3"))
$ws3.Range("B3").Value = "/room(GET)"
$ws3.Range("C3").Value = "add extra parameter for the filter as roomType in the server side"
$ws3.Range("D3").Value = "rajendra"
Set-DateStyle($ws3.Range("E3:F3"))
$ws3.Range("E3").Value = "19-9-16"
$ws3.Range("F3").Value = "19-9-16"
$ws3.Rows.Item(3).RowHeight = 30

# Row 4
Set-DateStyle($ws3.Range("A4"))
$ws3.Range("A4").Value = 42633
Set-TextStyle($ws3.Range("B4:D4"))
$ws3.Range("B4").Value = "/room(POST)"
$ws3.Range("C4").Value = "while saving room make the room type as mandaroty to select any one , in client side"
$ws3.Range("D4").Value = "rajashree"
Set-DateStyle($ws3.Range("E4:F4"))
$ws3.Range("E4").Value = "19-9-16"
$ws3.Range("F4").Value = "19-9-16"
$ws3.Rows.Item(4).RowHeight = 30

# Row 5
Set-DateStyle($ws3.Range("A5"))
$ws3.Range("A5").Value = 42634
Set-TextStyle($ws3.Range("B5:D5"))
$ws3.Range("B5").Value = "/room(POST)"
$ws3.Range("C5").Value = "while saving room make the room type as mandaroty to select any one , in server side"
$ws3.Range("D5").Value = "rajendra"
Set-DateStyle($ws3.Range("E5:F5"))
$ws3.Range("E5").Value = "19-9-16"
$ws3.Range("F5").Value = "19-9-16"
$ws3.Rows.Item(5).RowHeight = 30

# Row 6
Set-DateStyle($ws3.Range("A6"))
$ws3.Range("A6").Value = 42635
Set-TextStyle($ws3.Range("C6:D6"))
$ws3.Range("C6").Value = "In the client side show the validation error so that the user could know where the validation is not proper"
$ws3.Range("D6").Value = "rajashree"
Set-DateStyle($ws3.Range("E6:F6"))
$ws3.Range("E6").Value = "19-9-16"
$ws3.Range("F6").Value = "19-9-16"
Set-TextStyle($ws3.Range("H6"))
$ws3.Range("H6").Value = "This can be done using the alert of the ui bootstrap"
$ws3.Rows.Item(6).RowHeight = 45

# Row 7
Set-DateStyle($ws3.Range("A7"))
$ws3.Range("A7").Value = 42636
Set-TextStyle($ws3.Range("B7"))
$ws3.Range("B7").Value = "Error not proper for the auth in the server console"
Set-TextStyle($ws3.Range("D7"))
$ws3.Range("D7").Value = "rajendra"
Set-DateStyle($ws3.Range("E7:F7"))
$ws3.Range("E7").Value = "19-9-16"
$ws3.Range("F7").Value = "19-9-16"
$ws3.Rows.Item(7).RowHeight = 45

# Rows 8-20: just the running date series in column A.
$dateSerials = @(42637, 42638, 42639, 42640, 42641, 42642, 42643, 42644, 42645, 42646, 42647, 42648, 42649)
$row = 8
foreach ($serial in $dateSerials) {
    Set-DateStyle($ws3.Cells.Item($row, 1))
    $ws3.Cells.Item($row, 1).Value = $serial
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 2. Selection / active-sheet bookkeeping to match the new layout.
# ---------------------------------------------------------------------------
$ws2.Cells.Select()
$ws3.Range("D7:F7").Select()
$ws3.Activate()
